$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(20).Delete()
